# Auto-generated edit script to apply cell value changes per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("O2").Value2 = 1.22
$ws.Range("Q2").Value2 = 1.26
$ws.Range("S2").Value2 = 1.26
# Row 4
$ws.Range("G4").Value2 = 6.4
$ws.Range("Q4").Value2 = 2.46
$ws.Range("W4").Value2 = 1.19
# Row 5
$ws.Range("AH5").Value2 = 980
$ws.Range("AI5").Value2 = 980
# Row 6
$ws.Range("F6").Value2 = 4
$ws.Range("X6").Value2 = 1000
$ws.Range("Y6").Value2 = 1000
$ws.Range("Z6").Value2 = 1000
$ws.Range("AB6").Value2 = 1000
$ws.Range("AD6").Value2 = 1000
$ws.Range("AO6").Value2 = 1000
# Row 8
$ws.Range("F8").Value2 = 1.95
$ws.Range("Y8").Value2 = 1000
$ws.Range("AI8").Value2 = 65
# Row 9
$ws.Range("H9").Value2 = 2.18
$ws.Range("K9").Value2 = 3.8
$ws.Range("Q9").Value2 = 1.84
# Row 10
$ws.Range("F10").Value2 = 3.75
$ws.Range("G10").Value2 = 5.2
$ws.Range("I10").Value2 = 2.06
$ws.Range("K10").Value2 = 4.7
$ws.Range("Q10").Value2 = 1.53
$ws.Range("V10").Value2 = 1.94
$ws.Range("W10").Value2 = 1.24
$ws.Range("Y10").Value2 = 1000
$ws.Range("Z10").Value2 = 1000
$ws.Range("AC10").Value2 = 1000
$ws.Range("AD10").Value2 = 1000
# Row 11
$ws.Range("G11").Value2 = 1.3
$ws.Range("K11").Value2 = 7.4
$ws.Range("R11").Value2 = 1.51
$ws.Range("S11").Value2 = 1.53
$ws.Range("T11").Value2 = 1.74
$ws.Range("W11").Value2 = 4.2
# Row 13
$ws.Range("F13").Value2 = 2.22
$ws.Range("G13").Value2 = 2.68
$ws.Range("I13").Value2 = 3.6
$ws.Range("N13").Value2 = 4.8
$ws.Range("R13").Value2 = 1.56
$ws.Range("V13").Value2 = 1.39
$ws.Range("W13").Value2 = 1.6
# Row 14
$ws.Range("F14").Value2 = 10.5
$ws.Range("G14").Value2 = 11.5
$ws.Range("H14").Value2 = 1.38
$ws.Range("I14").Value2 = 1.39
$ws.Range("J14").Value2 = 5.2
$ws.Range("K14").Value2 = 5.4
# Row 15
$ws.Range("G15").Value2 = 5.5
$ws.Range("H15").Value2 = 1.71
$ws.Range("K15").Value2 = 5.1
$ws.Range("L15").Value2 = 1.21
# Row 17
$ws.Range("S17").Value2 = 2.18
$ws.Range("T17").Value2 = 1.74
$ws.Range("U17").Value2 = 1.52
# Row 18
$ws.Range("F18").Value2 = 2.56
$ws.Range("G18").Value2 = 2.58
$ws.Range("N18").Value2 = 3.45
$ws.Range("O18").Value2 = 1.38
$ws.Range("P18").Value2 = 1.84
$ws.Range("Q18").Value2 = 2.16
$ws.Range("W18").Value2 = 1.63
$ws.Range("AE18").Value2 = 36
$ws.Range("AK18").Value2 = 29
$ws.Range("AO18").Value2 = 36
# Row 19
$ws.Range("F19").Value2 = 1.87
$ws.Range("G19").Value2 = 1.88
$ws.Range("I19").Value2 = 5.4
$ws.Range("L19").Value2 = 1.01
$ws.Range("N19").Value2 = 3
$ws.Range("T19").Value2 = 2.18
$ws.Range("X19").Value2 = 10.5
$ws.Range("AA19").Value2 = 160
$ws.Range("AD19").Value2 = 22
$ws.Range("AI19").Value2 = 110
$ws.Range("AN19").Value2 = 17.5
# Row 20
$ws.Range("F20").Value2 = 1.34
$ws.Range("I20").Value2 = 12
$ws.Range("L20").Value2 = 1.34
$ws.Range("M20").Value2 = 1.01
$ws.Range("N20").Value2 = 3.9
$ws.Range("O20").Value2 = 1.24
$ws.Range("R20").Value2 = 1.39
$ws.Range("S20").Value2 = 2.6
$ws.Range("T20").Value2 = 2
$ws.Range("U20").Value2 = 1.62
$ws.Range("V20").Value2 = 1.09
$ws.Range("W20").Value2 = 3.85
$ws.Range("X20").Value2 = 21
$ws.Range("Y20").Value2 = 38
$ws.Range("Z20").Value2 = 120
$ws.Range("AA20").Value2 = 590
$ws.Range("AB20").Value2 = 9.199999999999999
$ws.Range("AC20").Value2 = 14.5
$ws.Range("AD20").Value2 = 46
$ws.Range("AE20").Value2 = 240
$ws.Range("AF20").Value2 = 8.800000000000001
$ws.Range("AG20").Value2 = 11.5
$ws.Range("AH20").Value2 = 36
$ws.Range("AI20").Value2 = 200
$ws.Range("AJ20").Value2 = 11.5
$ws.Range("AK20").Value2 = 17
$ws.Range("AL20").Value2 = 46
$ws.Range("AM20").Value2 = 230
$ws.Range("AN20").Value2 = 6.6
$ws.Range("AO20").Value2 = 310

$wb.Save()
